$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: S3 changes from 1900 to 2000
$ws.Range("S3").Value = 2000

# Rows 32-74 (wide scan): set Obj-posXpix (R), Obj-posYpix (S), run (W), quality (X)
$ws.Range("R32").Value = 2700
$ws.Range("S32").Value = 3400
$ws.Range("W32").Value = 1
$ws.Range("X32").Value = 1
$ws.Range("R33").Value = 2100
$ws.Range("S33").Value = 200
$ws.Range("W33").Value = 1
$ws.Range("X33").Value = 1
$ws.Range("R34").Value = 2100
$ws.Range("S34").Value = 800
$ws.Range("W34").Value = 1
$ws.Range("X34").Value = 1
$ws.Range("R35").Value = 2100
$ws.Range("S35").Value = 1400
$ws.Range("W35").Value = 1
$ws.Range("X35").Value = 1
$ws.Range("R36").Value = 2200
$ws.Range("S36").Value = 2000
$ws.Range("W36").Value = 1
$ws.Range("X36").Value = 1
$ws.Range("R37").Value = 2200
$ws.Range("S37").Value = 2700
$ws.Range("W37").Value = 1
$ws.Range("X37").Value = 1
$ws.Range("R38").Value = 2200
$ws.Range("S38").Value = 3300
$ws.Range("W38").Value = 1
$ws.Range("X38").Value = 1
$ws.Range("R39").Value = 1500
$ws.Range("S39").Value = 100
$ws.Range("W39").Value = 1
$ws.Range("X39").Value = 1
$ws.Range("R40").Value = 1500
$ws.Range("S40").Value = 700
$ws.Range("W40").Value = 1
$ws.Range("X40").Value = 1
$ws.Range("R41").Value = 1500
$ws.Range("S41").Value = 1300
$ws.Range("W41").Value = 1
$ws.Range("X41").Value = 1
$ws.Range("R42").Value = 1500
$ws.Range("S42").Value = 1900
$ws.Range("W42").Value = 1
$ws.Range("X42").Value = 1
$ws.Range("R43").Value = 1500
$ws.Range("S43").Value = 2500
$ws.Range("W43").Value = 1
$ws.Range("X43").Value = 1
$ws.Range("R44").Value = 1500
$ws.Range("S44").Value = 3100
$ws.Range("W44").Value = 1
$ws.Range("X44").Value = 1
$ws.Range("R46").Value = 800
$ws.Range("S46").Value = 600
$ws.Range("W46").Value = 1
$ws.Range("X46").Value = 1
$ws.Range("R47").Value = 800
$ws.Range("S47").Value = 1200
$ws.Range("W47").Value = 1
$ws.Range("X47").Value = 1
$ws.Range("R48").Value = 800
$ws.Range("S48").Value = 1800
$ws.Range("W48").Value = 1
$ws.Range("X48").Value = 1
$ws.Range("R49").Value = 800
$ws.Range("S49").Value = 2400
$ws.Range("W49").Value = 1
$ws.Range("X49").Value = 1
$ws.Range("R50").Value = 800
$ws.Range("S50").Value = 3000
$ws.Range("W50").Value = 1
$ws.Range("X50").Value = 1
$ws.Range("R52").Value = 300
$ws.Range("S52").Value = 500
$ws.Range("W52").Value = 1
$ws.Range("X52").Value = 1
$ws.Range("R53").Value = 300
$ws.Range("S53").Value = 1100
$ws.Range("W53").Value = 1
$ws.Range("X53").Value = 1
$ws.Range("R54").Value = 200
$ws.Range("S54").Value = 1700
$ws.Range("W54").Value = 1
$ws.Range("X54").Value = 1
$ws.Range("R55").Value = 300
$ws.Range("S55").Value = 2300
$ws.Range("W55").Value = 1
$ws.Range("X55").Value = 1
$ws.Range("R56").Value = 200
$ws.Range("S56").Value = 2900
$ws.Range("W56").Value = 1
$ws.Range("X56").Value = 1
$ws.Range("R57").Value = 200
$ws.Range("S57").Value = 3500
$ws.Range("W57").Value = 1
$ws.Range("X57").Value = 1
$ws.Range("R58").Value = 300
$ws.Range("S58").Value = 2200
$ws.Range("W58").Value = 1
$ws.Range("X58").Value = 1
$ws.Range("R59").Value = 300
$ws.Range("S59").Value = 2300
$ws.Range("W59").Value = 1
$ws.Range("X59").Value = 1
$ws.Range("R62").Value = 700
$ws.Range("S62").Value = 100
$ws.Range("W62").Value = 1
$ws.Range("X62").Value = 1
$ws.Range("R63").Value = 700
$ws.Range("S63").Value = 200
$ws.Range("W63").Value = 1
$ws.Range("X63").Value = 1
$ws.Range("R64").Value = 700
$ws.Range("S64").Value = 300
$ws.Range("W64").Value = 1
$ws.Range("X64").Value = 1
$ws.Range("R65").Value = 700
$ws.Range("S65").Value = 400
$ws.Range("W65").Value = 1
$ws.Range("X65").Value = 1
$ws.Range("R66").Value = 700
$ws.Range("S66").Value = 500
$ws.Range("W66").Value = 1
$ws.Range("X66").Value = 1
$ws.Range("R70").Value = 600
$ws.Range("S70").Value = 100
$ws.Range("W70").Value = 1
$ws.Range("X70").Value = 1
$ws.Range("R71").Value = 600
$ws.Range("S71").Value = 200
$ws.Range("W71").Value = 1
$ws.Range("X71").Value = 1
$ws.Range("R72").Value = 650
$ws.Range("S72").Value = 250
$ws.Range("W72").Value = 1
$ws.Range("X72").Value = 1
$ws.Range("R73").Value = 650
$ws.Range("S73").Value = 350
$ws.Range("W73").Value = 1
$ws.Range("X73").Value = 1
$ws.Range("R74").Value = 1300
$ws.Range("S74").Value = 3400
$ws.Range("W74").Value = 1
$ws.Range("X74").Value = 1
